# The "Baker, Salas and Nelson Incubator" record (id=2) was removed from the
# incubators data set. Deleting its row shifts every subsequent record up by
# one row (so the old id=3 row becomes the new row 2, etc.) and shrinks the
# sheet's used range from A1:H30 down to A1:H29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the entire second row (the row right after the header row),
# shifting all the rows below it up by one.
$ws.Rows(2).Delete()

# Keep the "numbers stored as text" error-checking suppression in sync with
# the new (smaller) data range now that a row has been removed.
try {
    $ws.Range("A1:H29").Errors.Item(9).Ignore = $true
} catch {
}
